$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "...a young cellist named Jacques Offenbach (16)."
#    -> "...a young cellist named Jacob (Jacques) Offenbach (16)."
#    The replaced name is split out into its own run (same character
#    formatting as its neighbours: color 000000), matching how the phrase
#    is now made up of three separate runs in the edited document.
# ---------------------------------------------------------------------------

$target = $d.Content
$target.Find.Execute("Jacques Offenbach", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$nameStart = $target.Start
$nameEnd = $target.End

$nameRange = $d.Range($nameStart, $nameEnd)
$nameRange.Text = "Jacob (Jacques) Offenbach"

# Force the newly-typed name to live in its own run, distinct from the
# surrounding text (which keeps the original "color 000000" formatting).
$newNameLen = "Jacob (Jacques) Offenbach".Length
$nameRange2 = $d.Range($nameStart, $nameStart + $newNameLen)
$nameRange2.Bold = $true
$nameRange2.Bold = $false

# ---------------------------------------------------------------------------
# 2) Date stamp near the end of the document: "9 April 2016" -> "9 June 2016"
# ---------------------------------------------------------------------------

$afterCopyright = $d.Content
$afterCopyright.Find.Execute("Scharfenberger", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dateScope = $d.Range($afterCopyright.End, $d.Content.End)
$dateScope.Find.Execute("April", $false, $false, $false, $false, $false, $true, 1, $false, "June", 2)

# Re-locate the just-replaced month so it keeps its own run (distinct from
# the following " 2016" run) instead of being silently merged with it.
$monthScope = $d.Range($afterCopyright.End, $d.Content.End)
$monthScope.Find.Execute("June", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$monthRange = $d.Range($monthScope.Start, $monthScope.End)
$monthRange.Bold = $true
$monthRange.Bold = $false
